# Slide 3 ("18 - Descrição de processos de Negócios"): the "Evento" paragraph
# goes from
#     [bold]"Evento: "[/bold][normal]"Fábrica trata resposta ao orçamento"[/normal]
# to
#     [bold]"Evento"[/bold][bold]": "[/bold][normal]"Loja "[/normal][normal]"trata resposta ao orçamento"[/normal]
#
# We locate the shape/text dynamically (rather than hard-coding character
# offsets) and perform the edits back-to-front so earlier offsets stay valid.

$p = $ppt.ActivePresentation

# Find the slide + shape that holds the "Evento: Fábrica ..." text.
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $txt = $shape.TextFrame.TextRange.Text
            if ($txt -ne $null -and $txt.Contains("Fábrica trata resposta ao orçamento")) {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# --- Step 1: split "Fábrica trata resposta ao orçamento" into
#             "Loja " (no dirty, not bold) + "trata resposta ao orçamento" (dirty, not bold)
$full = $tr.Text
$oldTail = "Fábrica trata resposta ao orçamento"
$idx = $full.IndexOf($oldTail)
$startPos = $idx + 1                      # 1-based start of "Fábrica trata resposta ao orçamento"
$prefixLen = "Fábrica ".Length            # length of "Fábrica " (8)

# Replace just the "Fábrica " prefix with "Loja " -- this keeps the
# remaining "trata resposta ao orçamento" run/run-properties untouched.
$prefixRange = $tr.Characters($startPos, $prefixLen)
$prefixRange.Text = "Loja "

# --- Step 2: split "Evento: " into "Evento" (bold, unchanged) + ": " (bold)
$full = $tr.Text
$oldHead = "Evento: "
$idx2 = $full.IndexOf($oldHead)
$colonStart = $idx2 + 1 + "Evento".Length  # 1-based start of ": "
$colonLen = ": ".Length

$colonRange = $tr.Characters($colonStart, $colonLen)
# Re-assert bold explicitly so the split-off run carries b="1" like the
# "Evento" run it came from.
$colonRange.Font.Bold = $true
$colonRange.Text = ": "

Write-Host "Final paragraph text:" $tr.Text
